$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4316.174
$ws.Range("J116").Value = 3753.2307
$ws.Range("L116").Value = 3753.2307
$ws.Range("N116").Value = -10637.2307

$ws.Range("H125").Value = 1237.5
$ws.Range("I125").Value = 1750
$ws.Range("J125").Value = 725
$ws.Range("K125").Value = 15750
$ws.Range("L125").Value = 6525
$ws.Range("M125").Value = -13290
$ws.Range("N125").Value = -11445

$ws.Range("H137").Value = 958.75
$ws.Range("I137").Value = 866.6667
$ws.Range("J137").Value = 1235
$ws.Range("K137").Value = 2600.0001
$ws.Range("L137").Value = 3705
$ws.Range("M137").Value = -50.0001000000002
$ws.Range("N137").Value = -8805

$ws.Range("H138").Value = 1862.13
$ws.Range("I138").Value = 1134.7826
$ws.Range("J138").Value = 2079.3896
$ws.Range("K138").Value = 3404.3478
$ws.Range("L138").Value = 6238.168799999999
$ws.Range("M138").Value = 1735.6522
$ws.Range("N138").Value = -16518.1688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H125").Value = 48600
$ws.Range("J125").Value = 48600
$ws.Range("L125").Value = 48600
$ws.Range("N125").Value = -58440

$ws.Range("H134").Value = 29317.447
$ws.Range("I134").Value = 2854.8235
$ws.Range("J134").Value = 254249.75
$ws.Range("K134").Value = 8564.470499999999
$ws.Range("L134").Value = 762749.25
$ws.Range("M134").Value = -6029.470499999999
$ws.Range("N134").Value = -767819.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2519.6
$ws.Range("I3").Value = 2799
$ws.Range("J3").Value = 2333.3333
$ws.Range("K3").Value = 2799
$ws.Range("L3").Value = 2333.3333
$ws.Range("M3").Value = -2686
$ws.Range("N3").Value = -2559.3333

$ws.Range("H10").Value = 666.3333
$ws.Range("I10").Value = 749.5
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 749.5
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = -610.5
$ws.Range("N10").Value = -778

$ws.Range("H12").Value = 250
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H14").Value = 10000
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10340

$ws.Range("H19").Value = 548.3333
$ws.Range("I19").Value = 262
$ws.Range("K19").Value = 262
$ws.Range("M19").Value = -92
$ws.Range("N19").ClearContents()

$ws.Range("H24").Value = 548.3333
$ws.Range("I24").Value = 262
$ws.Range("K24").Value = 262
$ws.Range("M24").Value = -92
$ws.Range("N24").ClearContents()

$ws.Range("H132").Value = 1707.36
$ws.Range("I132").Value = 1076.3889
$ws.Range("K132").Value = 3229.1667
$ws.Range("M132").Value = -699.1666999999998
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1002.5
$ws.Range("I55").Value = 800
$ws.Range("J55").Value = 1034.4736
$ws.Range("K55").Value = 2400
$ws.Range("L55").Value = 3103.4208
$ws.Range("M55").Value = -2223
$ws.Range("N55").Value = -3457.4208

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 49497.43
$ws.Range("J123").Value = 49497.43
$ws.Range("L123").Value = 49497.43
$ws.Range("N123").Value = -54397.43

$ws.Range("H126").Value = 4521.3076
$ws.Range("I126").Value = 3876.3
$ws.Range("K126").Value = 11628.9
$ws.Range("M126").Value = -9158.900000000001
$ws.Range("N126").ClearContents()

$ws.Range("H131").Value = 27000
$ws.Range("J131").Value = 27000
$ws.Range("L131").Value = 27000
$ws.Range("N131").Value = -37080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 398.2857
$ws.Range("I9").Value = 572.25
$ws.Range("J9").Value = 166.33333
$ws.Range("K9").Value = 572.25
$ws.Range("L9").Value = 166.33333
$ws.Range("M9").Value = -348.25
$ws.Range("N9").Value = -614.3333299999999

$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H14").Value = 1335320
$ws.Range("I14").Value = 4000000
$ws.Range("J14").Value = 2980
$ws.Range("K14").Value = 4000000
$ws.Range("L14").Value = 2980
$ws.Range("M14").Value = -3999828
$ws.Range("N14").Value = -3324

$ws.Range("H17").Value = 1830
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2037.5
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 2037.5
$ws.Range("M17").Value = -830
$ws.Range("N17").Value = -2377.5

$ws.Range("H19").Value = 903
$ws.Range("I19").Value = 903
$ws.Range("K19").Value = 903
$ws.Range("M19").Value = -733

$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2812
$ws.Range("N46").ClearContents()

$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 20000
$ws.Range("N112").Value = -22954

$ws.Range("H122").Value = 3859.6
$ws.Range("I122").Value = 2324.5
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 6973.5
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -4523.5
$ws.Range("N122").Value = -34900

$ws.Range("H132").Value = 2071.5278
$ws.Range("I132").Value = 1282.8182
$ws.Range("J132").Value = 3310.9285
$ws.Range("K132").Value = 3848.4546
$ws.Range("L132").Value = 9932.7855
$ws.Range("M132").Value = -1318.4546
$ws.Range("N132").Value = -14992.7855

$ws.Range("H136").Value = 3593.5938
$ws.Range("I136").Value = 1865.2174
$ws.Range("J136").Value = 8010.5557
$ws.Range("K136").Value = 5595.6522
$ws.Range("L136").Value = 24031.6671
$ws.Range("M136").Value = -3045.6522
$ws.Range("N136").Value = -29131.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 849.4
$ws.Range("I126").Value = 866
$ws.Range("J126").Value = 700
$ws.Range("K126").Value = 2598
$ws.Range("L126").Value = 2100
$ws.Range("M126").Value = -128
$ws.Range("N126").Value = -7040

$ws.Range("H132").Value = 1266.8636
$ws.Range("I132").Value = 988.5333000000001
$ws.Range("K132").Value = 2965.5999
$ws.Range("M132").Value = -435.5999000000002
$ws.Range("N132").ClearContents()
